# Prefix the period/year header labels (row 1, columns B:E) on each sheet
# so Power BI can automatically promote the first row to headers.
#
# Sheets "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Emissoes Totais (MtCO2eq)" and
# "Custo Total (bilhões de R$)" use plain years -> prefix with "Ano ".
# Sheet "Potencia Incremental - SIN(MW)" uses year ranges -> prefix with
# "Intervalo ".

$wb = $excel.ActiveWorkbook

$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)",
    "Custo Total (bilhões de R$)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($col in @("B", "C", "D", "E")) {
        $cell = $ws.Range($col + "1")
        $val = $cell.Value2
        if ($val -ne $null -and $val -ne "") {
            $cell.Value = "Ano " + $val
        }
    }
}

$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
foreach ($col in @("B", "C", "D", "E")) {
    $cell = $ws4.Range($col + "1")
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = "Intervalo " + $val
    }
}

$wb.Save()
